# Apply the weekly crime-data refresh to the CompStat worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number and reporting week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  52"
$ws.Range("C9").Value = "Report Covering the Week  12/23/2024  Through  12/29/2024"

# --- Reference cells used to copy exact cell formatting (styles) ---
# Style 13 (text placeholder, e.g. "0"/"***.*"): reference C22
# Style 14 (whole-number count): reference C16
# Style 15 (percentage change): reference E16
$styleTextRef = $ws.Range("C22")
$styleNumRef = $ws.Range("C16")
$stylePctRef = $ws.Range("E16")

# --- Cells that change FROM a numeric value TO the text placeholder "0" ---
# (style must switch from the numeric style to the text style)
$ws.Range("C15").Value = "'0"
$ws.Range("F22").Value = "'0"
$ws.Range("C27").Value = "'0"
$ws.Range("C28").Value = "'0"
$ws.Range("D28").Value = "'0"
$ws.Range("F33").Value = "'0"
$styleTextRef.Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("F33").PasteSpecial(-4122)

# --- Cell that changes FROM a numeric percentage TO the text placeholder "***.*" ---
$ws.Range("E28").Value = "***.*"
$styleTextRef.Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122)

# --- Cells that change FROM the text placeholder "0" TO a numeric count ---
$ws.Range("D15").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("D31").Value = 2
$styleNumRef.Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D31").PasteSpecial(-4122)

# --- Cells that change FROM the text placeholder "***.*" TO a numeric percentage ---
$ws.Range("E15").Value = -100
$ws.Range("E27").Value = -100
$ws.Range("E31").Value = -100
$stylePctRef.Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E31").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Plain value updates (style unchanged) ---
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -66.666666666666
$ws.Range("J15").Value = 24
$ws.Range("K15").Value = -37.5
$ws.Range("C16").Value = 3
$ws.Range("E16").Value = 200
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = -28.571428571428
$ws.Range("I16").Value = 108
$ws.Range("J16").Value = 92
$ws.Range("K16").Value = 17.391304347826
$ws.Range("L16").Value = 14.893617021276
$ws.Range("M16").Value = -39.325842696629
$ws.Range("N16").Value = -85.483870967741
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = -27.777777777777
$ws.Range("I17").Value = 199
$ws.Range("J17").Value = 221
$ws.Range("K17").Value = -9.954751131221
$ws.Range("L17").Value = -7.441860465116
$ws.Range("M17").Value = 17.058823529411
$ws.Range("N17").Value = -36.624203821656
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -25
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -43.75
$ws.Range("I18").Value = 127
$ws.Range("J18").Value = 121
$ws.Range("K18").Value = 4.958677685950
$ws.Range("L18").Value = -32.085561497326
$ws.Range("M18").Value = -66.840731070496
$ws.Range("N18").Value = -92.467378410438
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -30
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 39
$ws.Range("H19").Value = -17.948717948717
$ws.Range("I19").Value = 524
$ws.Range("J19").Value = 551
$ws.Range("K19").Value = -4.900181488203
$ws.Range("L19").Value = -18.125
$ws.Range("M19").Value = 35.051546391752
$ws.Range("N19").Value = -20.485584218512
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -25
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = -35.714285714285
$ws.Range("I20").Value = 188
$ws.Range("J20").Value = 156
$ws.Range("K20").Value = 20.512820512820
$ws.Range("L20").Value = 64.912280701754
$ws.Range("M20").Value = 39.259259259259
$ws.Range("N20").Value = -89.624724061810
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = -12.5
$ws.Range("G21").Value = 97
$ws.Range("H21").Value = -28.865979381443
$ws.Range("I21").Value = 1166
$ws.Range("J21").Value = 1167
$ws.Range("K21").Value = -0.085689802913
$ws.Range("L21").Value = -8.261211644374
$ws.Range("M21").Value = -7.826086956521
$ws.Range("N21").Value = -77.765064836003
$ws.Range("L22").Value = 166.666666666667
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = -10.526315789473
$ws.Range("F24").Value = 74
$ws.Range("G24").Value = 67
$ws.Range("H24").Value = 10.447761194029
$ws.Range("I24").Value = 1011
$ws.Range("J24").Value = 1048
$ws.Range("K24").Value = -3.530534351145
$ws.Range("L24").Value = -8.340888485947
$ws.Range("M24").Value = 13.595505617977
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = -33.333333333333
$ws.Range("G25").Value = 15
$ws.Range("H25").Value = -13.333333333333
$ws.Range("I25").Value = 227
$ws.Range("J25").Value = 311
$ws.Range("K25").Value = -27.009646302250
$ws.Range("L25").Value = -31.419939577039
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 9
$ws.Range("F26").Value = 24
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = -20
$ws.Range("I26").Value = 429
$ws.Range("J26").Value = 351
$ws.Range("K26").Value = 22.222222222222
$ws.Range("L26").Value = 26.176470588235
$ws.Range("M26").Value = 7.25
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -80
$ws.Range("J27").Value = 31
$ws.Range("K27").Value = -19.354838709677
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = -50
$ws.Range("L28").Value = -16.25
$ws.Range("G31").Value = 3
$ws.Range("J31").Value = 13
$ws.Range("K31").Value = 38.461538461538
$ws.Range("L33").Value = 25
